$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Données du mois d'août (rows 104-134), colonnes A (Date), B (CA), C (ACHATS)
$data = @(
    @(45139, 516000, 43000),
    @(45140, 2514750, 823500),
    @(45141, 1210000, 408150),
    @(45142, 1610750, 478000),
    @(45143, 1205750, 452650),
    @(45144, 1454750, 2012300),
    @(45145, 582750, 435450),
    @(45146, 602500, 123900),
    @(45147, 2105750, 764500),
    @(45148, 1129500, 1068450),
    @(45149, 1368250, 556600),
    @(45150, 2536500, 433100),
    @(45151, 1090500, 2517275),
    @(45152, 1984500, 246775),
    @(45153, 457250, 42800),
    @(45154, 2845600, 1334150),
    @(45155, 1144500, 704300),
    @(45156, 1302500, 184700),
    @(45157, 1354500, 118050),
    @(45158, 34500, 2252525),
    @(45159, 663500, 1095000),
    @(45160, 588500, 265400),
    @(45161, 2563250, 1081100),
    @(45162, 1214000, 374250),
    @(45163, 1177750, 218900),
    @(45164, 1987250, 780025),
    @(45165, 1166750, 1892300),
    @(45166, 1240750, 607000),
    @(45167, 787500, 106000),
    @(45168, 2303250, 1691750),
    @(45169, 950500, 427450),
)

$startRow = 104
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
